# Update "想去人数" (Column F) counts on the "展览" and "全部类型" sheets
# to reflect the latest generated data (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 4294
$ws1.Range("F3").Value = 103
$ws1.Range("F5").Value = 43
$ws1.Range("F7").Value = 42
$ws1.Range("F9").Value = 128
$ws1.Range("F10").Value = 311
$ws1.Range("F11").Value = 243
$ws1.Range("F12").Value = 2938
$ws1.Range("F13").Value = 143
$ws1.Range("F14").Value = 1505

# --- Sheet "全部类型" (sheet4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 4294
$ws4.Range("F3").Value = 103
$ws4.Range("F5").Value = 43
$ws4.Range("F8").Value = 42
$ws4.Range("F10").Value = 128
$ws4.Range("F11").Value = 311
$ws4.Range("F12").Value = 243
$ws4.Range("F13").Value = 2938
$ws4.Range("F14").Value = 143
$ws4.Range("F15").Value = 1505
